$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.032.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.948.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0853"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.410.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +51.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.945.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.021.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0957"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "266.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.164"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.10%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.110"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "33.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0445"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.26%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.116"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.41%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.88%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.18%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.272"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.997.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0330"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.47%  "
